{"js": "// The document contains a single hyperlink run reading \"4 Replies\"\n// (the count of comments on the blog post). The author bumped this to\n// \"5 Replies\" when the chapter's docx/pdf ebook content was regenerated.\n// Locate that exact text and update it in place, leaving all\n// surrounding formatting / the hyperlink itself untouched.\nconst results = context.document.body.search(\"4 Replies\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find the text \"4 Replies\" in the document.');\n}\n\nfor (const found of results.items) {\n  found.insertText(\"5 Replies\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single hyperlink run reading \"4 Replies\"\n# (the comment count on the blog post). The author bumped this to\n# \"5 Replies\" when the chapter's docx/pdf ebook content was regenerated.\n# Locate that exact text and update it in place, leaving all\n# surrounding formatting / the hyperlink itself untouched.\n$d = $word.ActiveDocument\n\n$range = $d.Content\n$find = $range.Find\n$find.ClearFormatting()\n$find.Text = \"4 Replies\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $range.Text = \"5 Replies\"\n}\n"}
